$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new annotation cells (disorder markers) in column C, and the
# "Pierson syndrome" summary note in column E for row 22.
$ws.Range("C17").Value = "(?)"
$ws.Range("C22").Value = "(?)"
$ws.Range("E22").Value = "(Pierson syndrome)"

# Re-create the "duplicate values" conditional formatting rule on column A
# so the workbook gets a fresh dxf entry (leaving the old, now-unused dxf
# in place) and the rule points at the new dxf.
$rng = $ws.Range("A1:A25")
$fcs = $rng.FormatConditions
$oldFc = $fcs.Item(1)
$color = $oldFc.Interior.Color
$oldFc.Delete()

$newFc = $fcs.AddUniqueValues()
$newFc.DupeUnique = 1
$newFc.Interior.Color = $color
$newFc.SetFirstPriority()

# Move the active selection to E17 (cosmetic, matches author's last
# selected cell when saving).
$ws.Range("E17").Select()
